$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 61.27353666666666
$ws.Range("H2").Value = 183.82061
$ws.Range("I2").Value = 0.3474604587406809
$ws.Range("J2").Value = 0.3474604587406808
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 101.6208433333333
$ws.Range("N2").Value = 304.86253
$ws.Range("O2").Value = 0.7748298251610983
$ws.Range("P2").Value = 0.7748298251610983
$ws.Range("Q2").Value = 6226.668470082588
$ws.Range("R2").Value = 56040.0162307433
$ws.Range("S2").Value = 0.2692227264964368
$ws.Range("T2").Value = 0.2692227264964367

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 61.27353666666666
$ws.Range("H3").Value = 183.82061
$ws.Range("I3").Value = 0.3474604587406809
$ws.Range("J3").Value = 0.3474604587406808
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 23.70024566666666
$ws.Range("N3").Value = 71.100737
$ws.Range("O3").Value = 0.1807075852140151
$ws.Range("P3").Value = 0.1807075852140151
$ws.Range("Q3").Value = 1452.197871865508
$ws.Range("R3").Value = 13069.78084678957
$ws.Range("S3").Value = 0.06278874045638236
$ws.Range("T3").Value = 0.06278874045638236

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 61.27353666666666
$ws.Range("H4").Value = 183.82061
$ws.Range("I4").Value = 0.3474604587406809
$ws.Range("J4").Value = 0.3474604587406808
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 5.831378333333333
$ws.Range("N4").Value = 17.494135
$ws.Range("O4").Value = 0.04446258962488651
$ws.Range("P4").Value = 0.04446258962488651
$ws.Range("Q4").Value = 357.3091741247055
$ws.Range("R4").Value = 3215.78256712235
$ws.Range("S4").Value = 0.0154489917878617
$ws.Range("T4").Value = 0.0154489917878617

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 89.66709900000001
$ws.Range("H5").Value = 269.001297
$ws.Range("I5").Value = 0.5084702637939138
$ws.Range("J5").Value = 0.5084702637939138
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 101.6208433333333
$ws.Range("N5").Value = 304.86253
$ws.Range("O5").Value = 0.7748298251610983
$ws.Range("P5").Value = 0.7748298251610983
$ws.Range("Q5").Value = 9112.04621963349
$ws.Range("R5").Value = 82008.41597670142
$ws.Range("S5").Value = 0.3939779255950558
$ws.Range("T5").Value = 0.3939779255950558

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 89.66709900000001
$ws.Range("H6").Value = 269.001297
$ws.Range("I6").Value = 0.5084702637939138
$ws.Range("J6").Value = 0.5084702637939138
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 23.70024566666666
$ws.Range("N6").Value = 71.100737
$ws.Range("O6").Value = 0.1807075852140151
$ws.Range("P6").Value = 0.1807075852140151
$ws.Range("Q6").Value = 2125.132274517321
$ws.Range("R6").Value = 19126.19047065589
$ws.Range("S6").Value = 0.09188443352333141
$ws.Range("T6").Value = 0.09188443352333141

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 89.66709900000001
$ws.Range("H7").Value = 269.001297
$ws.Range("I7").Value = 0.5084702637939138
$ws.Range("J7").Value = 0.5084702637939138
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 5.831378333333333
$ws.Range("N7").Value = 17.494135
$ws.Range("O7").Value = 0.04446258962488651
$ws.Range("P7").Value = 0.04446258962488651
$ws.Range("Q7").Value = 522.882778321455
$ws.Range("R7").Value = 4705.945004893096
$ws.Range("S7").Value = 0.02260790467552658
$ws.Range("T7").Value = 0.02260790467552658

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 25.40615466666667
$ws.Range("H8").Value = 76.218464
$ws.Range("I8").Value = 0.1440692774654054
$ws.Range("J8").Value = 0.1440692774654053
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 101.6208433333333
$ws.Range("N8").Value = 304.86253
$ws.Range("O8").Value = 0.7748298251610983
$ws.Range("P8").Value = 0.7748298251610983
$ws.Range("Q8").Value = 2581.794863083769
$ws.Range("R8").Value = 23236.15376775392
$ws.Range("S8").Value = 0.1116291730696058
$ws.Range("T8").Value = 0.1116291730696058

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 25.40615466666667
$ws.Range("H9").Value = 76.218464
$ws.Range("I9").Value = 0.1440692774654054
$ws.Range("J9").Value = 0.1440692774654053
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 23.70024566666666
$ws.Range("N9").Value = 71.100737
$ws.Range("O9").Value = 0.1807075852140151
$ws.Range("P9").Value = 0.1807075852140151
$ws.Range("Q9").Value = 602.1321070453297
$ws.Range("R9").Value = 5419.188963407967
$ws.Range("S9").Value = 0.02603441123430133
$ws.Range("T9").Value = 0.02603441123430132

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 25.40615466666667
$ws.Range("H10").Value = 76.218464
$ws.Range("I10").Value = 0.1440692774654054
$ws.Range("J10").Value = 0.1440692774654053
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 5.831378333333333
$ws.Range("N10").Value = 17.494135
$ws.Range("O10").Value = 0.04446258962488651
$ws.Range("P10").Value = 0.04446258962488651
$ws.Range("Q10").Value = 148.1528998565155
$ws.Range("R10").Value = 1333.37609870864
$ws.Range("S10").Value = 0.006405693161498228
$ws.Range("T10").Value = 0.006405693161498228
